# Methods getNumberOfRows, getNumberOfCells, LoadData and RunTest Completed
#
# - Add a new "vOutData" output column (G) used by the test-runner helpers.
# - Add two more test-case rows (CT 04 / CT 05).
# - Fix a data entry in the existing CT 03 row (RunTest -> Yes).
# - Remove the stale AutoFilter so the growing table isn't clipped to A1:F4.
# - Re-apply a simple header/grid look across the (now bigger) table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New header cell for column G
# ---------------------------------------------------------------------------
$ws.Cells.Item(1, 7).Value = "vOutData"

# ---------------------------------------------------------------------------
# 2. Fill in column G for the existing rows (2-4)
# ---------------------------------------------------------------------------
$ws.Cells.Item(2, 7).Formula = "=DATE(2020,4,13)"
$ws.Cells.Item(3, 7).Value = 2
# row 4 col G intentionally left blank

# ---------------------------------------------------------------------------
# 3. Correct existing data - CT 03 / RunTest should be "Yes"
# ---------------------------------------------------------------------------
$ws.Cells.Item(4, 2).Value = "Yes"

# ---------------------------------------------------------------------------
# 4. Two new rows: CT 04 and CT 05
# ---------------------------------------------------------------------------
$ws.Cells.Item(5, 1).Value = "CT 04"
$ws.Cells.Item(5, 2).Value = "No"
$ws.Cells.Item(5, 3).Value = "Jordana"
$ws.Cells.Item(5, 4).Value = "Jordana"
$ws.Cells.Item(5, 5).Value = "Lelles Moreira"
$ws.Cells.Item(5, 6).Value = "Jordana"
$ws.Cells.Item(5, 7).Value = 4

$ws.Cells.Item(6, 1).Value = "CT 05"
$ws.Cells.Item(6, 2).Value = "Yes"
$ws.Cells.Item(6, 3).Value = "Carlos Anthony"
$ws.Cells.Item(6, 4).Value = "Carlos Anthony"
$ws.Cells.Item(6, 5).Value = "Lelles Moreira"
$ws.Cells.Item(6, 6).Value = "Carlos Anthony"
$ws.Cells.Item(6, 7).Value = "aqui"

# ---------------------------------------------------------------------------
# 5. Drop the AutoFilter (table outgrew A1:F4) and fix up the now-dangling
#    _FilterDatabase defined name that Excel leaves pointing at the old
#    range -- it becomes #REF! once the filter is gone.
# ---------------------------------------------------------------------------
$ws.AutoFilterMode = $false
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Plan1!#REF!"
    }
}

# ---------------------------------------------------------------------------
# 6. Column G should match column F's width, and the used range is now
#    A1:G6.
# ---------------------------------------------------------------------------
$ws.Columns.Item(7).ColumnWidth = $ws.Columns.Item(6).ColumnWidth

# ---------------------------------------------------------------------------
# 7. Re-style the whole table with a simple, uniform look:
#    - thin grid everywhere
#    - header row: bold white-on-blue fill, medium top border, no bottom
#      border, wrap text
#    - data rows: normal font, wrap text
# ---------------------------------------------------------------------------
$fullRange = $ws.Range("A1:G6")
$fullRange.Borders.LineStyle = 0

$fullRange.WrapText = $true
$fullRange.Borders.LineStyle = 1
$fullRange.Borders.Weight = 2

$headerRange = $ws.Range("A1:G1")
$headerRange.Font.Bold = $true
$headerRange.Interior.Pattern = 1
$headerRange.Interior.ThemeColor = 3
$headerRange.Interior.TintAndShade = 0.6
$headerRange.Borders.Item(8).LineStyle = 1
$headerRange.Borders.Item(8).Weight = -4138
$headerRange.Borders.Item(9).LineStyle = 0

# Right edge of the old last column (F) kept its original "outer edge"
# medium border even after G was appended -- match that quirk, and give
# the new true last column (G) the same medium right edge.
$ws.Range("F1").Borders.Item(10).Weight = -4138
$ws.Range("G1").Borders.Item(10).Weight = -4138

$dataRange = $ws.Range("A2:G6")
$dataRange.Font.Bold = $false

# numeric display for the new date column (built-in date format #14)
$ws.Cells.Item(2, 7).NumberFormat = "mm-dd-yy"

# ---------------------------------------------------------------------------
# 8. Row heights were carrying an explicit "thick bottom border" row flag
#    from the old header/footer rows -- AutoFit clears that stale metadata
#    now that the border scheme has changed.
# ---------------------------------------------------------------------------
$ws.Rows.Item(1).AutoFit()
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(3).AutoFit()
$ws.Rows.Item(4).AutoFit()
$ws.Rows.Item(5).AutoFit()
$ws.Rows.Item(6).AutoFit()

# ---------------------------------------------------------------------------
# 9. Selection ends up on G14 (matches the recorded cursor position after
#    the edits were made).
# ---------------------------------------------------------------------------
$ws.Range("G14").Select()
